$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 210.17857
$ws.Range("I33").Value = 210.55556
$ws.Range("K33").Value = 210.55556
$ws.Range("M33").Value = 18.44443999999999

$ws.Range("H62").Value = 4062.2632
$ws.Range("I62").Value = 2466.6428
$ws.Range("K62").Value = 2466.6428
$ws.Range("M62").Value = -1842.6428

$ws.Range("H65").Value = 4062.2632
$ws.Range("I65").Value = 2466.6428
$ws.Range("K65").Value = 12333.214
$ws.Range("M65").Value = -9213.214

$ws.Range("H137").Value = 45457680
$ws.Range("I137").Value = 125001750
$ws.Range("K137").Value = 375005250
$ws.Range("M137").Value = -375002700

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1228.25
$ws.Range("I2").Value = 1249.091
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 1249.091
$ws.Range("L2").Value = 999
$ws.Range("M2").Value = -1136.091
$ws.Range("N2").Value = -1225

$ws.Range("H32").Value = 5278.6064
$ws.Range("I32").Value = 5494.4727
$ws.Range("J32").Value = 3299.8333
$ws.Range("K32").Value = 5494.4727
$ws.Range("L32").Value = 3299.8333
$ws.Range("M32").Value = -5207.4727
$ws.Range("N32").Value = -3873.8333

$ws.Range("H45").Value = 4763.357
$ws.Range("I45").Value = 2731
$ws.Range("J45").Value = 5892.4443
$ws.Range("K45").Value = 2731
$ws.Range("L45").Value = 5892.4443
$ws.Range("M45").Value = -2354
$ws.Range("N45").Value = -6646.4443

$ws.Range("H61").Value = 20591438
$ws.Range("I61").Value = 29169630
$ws.Range("J61").Value = 3775.1
$ws.Range("K61").Value = 29169630
$ws.Range("L61").Value = 3775.1
$ws.Range("M61").Value = -29169418
$ws.Range("N61").Value = -4199.1

$ws.Range("H116").Value = 1228.25
$ws.Range("I116").Value = 1249.091
$ws.Range("J116").Value = 999
$ws.Range("K116").Value = 1249.091
$ws.Range("L116").Value = 999
$ws.Range("M116").Value = 1044.909
$ws.Range("N116").Value = -5587

$ws.Range("H136").Value = 20591438
$ws.Range("I136").Value = 29169630
$ws.Range("J136").Value = 3775.1
$ws.Range("K136").Value = 87508890
$ws.Range("L136").Value = 11325.3
$ws.Range("M136").Value = -87506340
$ws.Range("N136").Value = -16425.3

$ws.Range("H140").Value = 200000
$ws.Range("I140").Value = 200000
$ws.Range("K140").Value = 200000
$ws.Range("M140").Value = -194820

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1228.25
$ws.Range("I3").Value = 1249.091
$ws.Range("J3").Value = 999
$ws.Range("K3").Value = 1249.091
$ws.Range("L3").Value = 999
$ws.Range("M3").Value = -1135.091
$ws.Range("N3").Value = -1227

$ws.Range("H99").Value = 1764.2354
$ws.Range("I99").Value = 1850.75
$ws.Range("J99").Value = 380
$ws.Range("K99").Value = 1850.75
$ws.Range("L99").Value = 380
$ws.Range("M99").Value = -352.75
$ws.Range("N99").Value = -3376

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22730934
$ws.Range("I31").Value = 38464864
$ws.Range("K31").Value = 38464864
$ws.Range("M31").Value = -38464569

$ws.Range("H34").Value = 22730934
$ws.Range("I34").Value = 38464864
$ws.Range("K34").Value = 38464864
$ws.Range("M34").Value = -38464662

$ws.Range("H62").Value = 13896669
$ws.Range("I62").Value = 6439.0557
$ws.Range("K62").Value = 6439.0557
$ws.Range("M62").Value = -5815.0557

$ws.Range("H65").Value = 13896669
$ws.Range("I65").Value = 6439.0557
$ws.Range("K65").Value = 32195.2785
$ws.Range("M65").Value = -29075.2785

$ws.Range("H122").Value = 4055.5557
$ws.Range("I122").Value = 7083.3335
$ws.Range("K122").Value = 21250.0005
$ws.Range("M122").Value = -18800.0005

$ws.Range("H129").Value = 53074.5
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 53074.5
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 53074.5
$ws.Range("N129").Value = -63074.5

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0

$ws.Range("H131").Value = 84274.664
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 84274.664
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 84274.664
$ws.Range("N131").Value = -94354.664

$ws.Range("H132").Value = 1870.6111
$ws.Range("I132").Value = 1286.4482
$ws.Range("J132").Value = 4290.7144
$ws.Range("K132").Value = 3859.3446
$ws.Range("L132").Value = 12872.1432
$ws.Range("M132").Value = -1329.3446
$ws.Range("N132").Value = -17932.1432

$ws.Range("H133").Value = 64750
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 64750
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 64750
$ws.Range("N133").Value = -69810

$ws.Range("H134").Value = 1735.3478
$ws.Range("I134").Value = 1710.238
$ws.Range("J134").Value = 1999
$ws.Range("K134").Value = 5130.714
$ws.Range("L134").Value = 5997
$ws.Range("M134").Value = -2595.714
$ws.Range("N134").Value = -11067

$ws.Range("H135").Value = 103399.4
$ws.Range("I135").Value = 65000
$ws.Range("J135").Value = 128999
$ws.Range("K135").Value = 65000
$ws.Range("L135").Value = 128999
$ws.Range("M135").Value = -59930
$ws.Range("N135").Value = -139139

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0

$ws.Range("H138").Value = 105000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 105000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 105000
$ws.Range("N138").Value = -115280

$ws.Range("H139").Value = 95000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 95000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0

$ws.Range("H141").Value = 311972.2
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 311972.2
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 311972.2
$ws.Range("N141").Value = -322332.2

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5806.2666
$ws.Range("I40").Value = 5459.1
$ws.Range("K40").Value = 5459.1
$ws.Range("M40").Value = -5323.1

$ws.Range("H122").Value = 3834.422
$ws.Range("I122").Value = 3433.838
$ws.Range("J122").Value = 5687.125
$ws.Range("K122").Value = 10301.514
$ws.Range("L122").Value = 17061.375
$ws.Range("M122").Value = -7851.514000000001
$ws.Range("N122").Value = -21961.375

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 27000
$ws.Range("J63").Value = 27000
$ws.Range("L63").Value = 27000
$ws.Range("N63").Value = -28248

$ws.Range("H66").Value = 27000
$ws.Range("J66").Value = 27000
$ws.Range("L66").Value = 81000
$ws.Range("N66").Value = -87240

$ws.Range("H122").Value = 3033.4375
$ws.Range("I122").Value = 2912.3076
$ws.Range("K122").Value = 8736.9228
$ws.Range("M122").Value = -6286.9228
